$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 28.01347424679749
$ws.Cells.Item(2, 3).Value = 27.21930555597951
$ws.Cells.Item(2, 4).Value = 15.37063265131838
$ws.Cells.Item(2, 5).Value = 16.82971055238288
$ws.Cells.Item(2, 7).Value = 3.779004842666364
$ws.Cells.Item(2, 9).Value = 37.23724022780687
$ws.Cells.Item(2, 10).Value = 9.646773456234229
$ws.Cells.Item(2, 14).Value = 21.50631683333983

$ws.Cells.Item(3, 2).Value = 27.5584124960966
$ws.Cells.Item(3, 3).Value = 26.71750360611718
$ws.Cells.Item(3, 4).Value = 15.33445668237109
$ws.Cells.Item(3, 5).Value = 16.79589350237164
$ws.Cells.Item(3, 7).Value = 3.785325032427704
$ws.Cells.Item(3, 9).Value = 37.09264833351401
$ws.Cells.Item(3, 10).Value = 9.66683883639344
$ws.Cells.Item(3, 14).Value = 21.54308640729557

$ws.Cells.Item(4, 2).Value = 27.28551717651133
$ws.Cells.Item(4, 3).Value = 26.41488329764672
$ws.Cells.Item(4, 4).Value = 15.31630997733369
$ws.Cells.Item(4, 5).Value = 16.77959089353637
$ws.Cells.Item(4, 7).Value = 3.789394139793544
$ws.Cells.Item(4, 9).Value = 37.01291208349151
$ws.Cells.Item(4, 10).Value = 9.680900404898306
$ws.Cells.Item(4, 14).Value = 21.56764001558074

$ws.Cells.Item(5, 2).Value = 27.17609143358405
$ws.Cells.Item(5, 3).Value = 26.29310867897208
$ws.Cells.Item(5, 4).Value = 15.30993730720798
$ws.Cells.Item(5, 5).Value = 16.77406824498082
$ws.Cells.Item(5, 7).Value = 3.791099985843601
$ws.Cells.Item(5, 9).Value = 36.98269717352343
$ws.Cells.Item(5, 10).Value = 9.687067311503892
$ws.Cells.Item(5, 14).Value = 21.57814079653971

$ws.Cells.Item(6, 2).Value = 27.15803307031556
$ws.Cells.Item(6, 3).Value = 26.27298623819205
$ws.Cells.Item(6, 4).Value = 15.3089408512996
$ws.Cells.Item(6, 5).Value = 16.77321885722769
$ws.Cells.Item(6, 7).Value = 3.791386125801705
$ws.Cells.Item(6, 9).Value = 36.97781767637728
$ws.Cells.Item(6, 10).Value = 9.688117662709157
$ws.Cells.Item(6, 14).Value = 21.57991427755861

$ws.Cells.Item(7, 2).Value = 27.28403402367448
$ws.Cells.Item(7, 3).Value = 26.41323452735183
$ws.Cells.Item(7, 4).Value = 15.31621989495512
$ws.Cells.Item(7, 5).Value = 16.77951187716336
$ws.Cells.Item(7, 7).Value = 3.789416952149768
$ws.Cells.Item(7, 9).Value = 37.01249536667687
$ws.Cells.Item(7, 10).Value = 9.680981807315291
$ws.Cells.Item(7, 14).Value = 21.56777963137771

$ws.Cells.Item(8, 2).Value = 27.85530382744012
$ws.Cells.Item(8, 3).Value = 27.04524511264587
$ws.Cells.Item(8, 4).Value = 15.35731565305547
$ws.Cells.Item(8, 5).Value = 16.81712420451592
$ws.Cells.Item(8, 7).Value = 3.781145078548259
$ws.Cells.Item(8, 9).Value = 37.18551122677197
$ws.Cells.Item(8, 10).Value = 9.653329884615326
$ws.Cells.Item(8, 14).Value = 21.51858337615542

$ws.Cells.Item(9, 2).Value = 29.02082482279943
$ws.Cells.Item(9, 3).Value = 28.32083038802714
$ws.Cells.Item(9, 4).Value = 15.47012257289813
$ws.Cells.Item(9, 5).Value = 16.92628216577521
$ws.Cells.Item(9, 7).Value = 3.76640780166941
$ws.Cells.Item(9, 9).Value = 37.59625148730019
$ws.Cells.Item(9, 10).Value = 9.612973747395348
$ws.Cells.Item(9, 14).Value = 21.43789304218822

$ws.Cells.Item(10, 2).Value = 29.89605828100476
$ws.Cells.Item(10, 3).Value = 29.27033252927406
$ws.Cells.Item(10, 4).Value = 15.57253106775739
$ws.Cells.Item(10, 5).Value = 17.02800559133801
$ws.Cells.Item(10, 7).Value = 3.756468457761238
$ws.Cells.Item(10, 9).Value = 37.9410130148321
$ws.Cells.Item(10, 10).Value = 9.591850956959444
$ws.Cells.Item(10, 14).Value = 21.38836813052048

$ws.Cells.Item(11, 2).Value = 30.296516662271
$ws.Cells.Item(11, 3).Value = 29.70295957121034
$ws.Cells.Item(11, 4).Value = 15.62331572260976
$ws.Cells.Item(11, 5).Value = 17.07892316928681
$ws.Cells.Item(11, 7).Value = 3.75213607713767
$ws.Cells.Item(11, 9).Value = 38.10701168631855
$ws.Cells.Item(11, 10).Value = 9.584108156004737
$ws.Cells.Item(11, 14).Value = 21.36798724124925

$ws.Cells.Item(12, 2).Value = 30.44834450832702
$ws.Cells.Item(12, 3).Value = 29.86672395689087
$ws.Cells.Item(12, 4).Value = 15.64314488122653
$ws.Cells.Item(12, 5).Value = 17.09886743345237
$ws.Cells.Item(12, 7).Value = 3.750522424562452
$ws.Cells.Item(12, 9).Value = 38.17116959750845
$ws.Cells.Item(12, 10).Value = 9.581445660337874
$ws.Cells.Item(12, 14).Value = 21.36058105110234

$ws.Cells.Item(13, 2).Value = 30.41563985386317
$ws.Cells.Item(13, 3).Value = 29.83145963006077
$ws.Cells.Item(13, 4).Value = 15.63884781541271
$ws.Cells.Item(13, 5).Value = 17.09454268406396
$ws.Cells.Item(13, 7).Value = 3.75086876018438
$ws.Cells.Item(13, 9).Value = 38.1572946675168
$ws.Cells.Item(13, 10).Value = 9.58200707002157
$ws.Cells.Item(13, 14).Value = 21.36216220405525

$ws.Cells.Item(14, 2).Value = 30.30900480199546
$ws.Cells.Item(14, 3).Value = 29.71643471225646
$ws.Cells.Item(14, 4).Value = 15.62493512789285
$ws.Cells.Item(14, 5).Value = 17.08055074163551
$ws.Cells.Item(14, 7).Value = 3.75200278284751
$ws.Cells.Item(14, 9).Value = 38.11226410964496
$ws.Cells.Item(14, 10).Value = 9.583883702248773
$ws.Cells.Item(14, 14).Value = 21.36737166875711

$ws.Cells.Item(15, 2).Value = 30.24370722973332
$ws.Cells.Item(15, 3).Value = 29.64596585862256
$ws.Cells.Item(15, 4).Value = 15.6164909140309
$ws.Cells.Item(15, 5).Value = 17.07206644099681
$ws.Cells.Item(15, 7).Value = 3.752700903183608
$ws.Cells.Item(15, 9).Value = 38.08484996983663
$ws.Cells.Item(15, 10).Value = 9.585068330117515
$ws.Cells.Item(15, 14).Value = 21.37060327815421

$ws.Cells.Item(16, 2).Value = 29.8699210080192
$ws.Cells.Item(16, 3).Value = 29.24205930421117
$ws.Cells.Item(16, 4).Value = 15.56929619608566
$ws.Cells.Item(16, 5).Value = 17.02477105181415
$ws.Cells.Item(16, 7).Value = 3.756755371285729
$ws.Cells.Item(16, 9).Value = 37.93034737070266
$ws.Cells.Item(16, 10).Value = 9.592394620204262
$ws.Cells.Item(16, 14).Value = 21.38974352673064

$ws.Cells.Item(17, 2).Value = 29.64109936243022
$ws.Cells.Item(17, 3).Value = 28.9943361251765
$ws.Cells.Item(17, 4).Value = 15.5414153983546
$ws.Cells.Item(17, 5).Value = 16.99694297494131
$ws.Cells.Item(17, 7).Value = 3.759290896630566
$ws.Cells.Item(17, 9).Value = 37.83789886454632
$ws.Cells.Item(17, 10).Value = 9.597367866675814
$ws.Cells.Item(17, 14).Value = 21.40203756436816

$ws.Cells.Item(18, 2).Value = 29.50971333638248
$ws.Cells.Item(18, 3).Value = 28.85192689148415
$ws.Cells.Item(18, 4).Value = 15.52577468300193
$ws.Cells.Item(18, 5).Value = 16.98137438047241
$ws.Cells.Item(18, 7).Value = 3.76076707639284
$ws.Cells.Item(18, 9).Value = 37.78558893642454
$ws.Cells.Item(18, 10).Value = 9.600403919467514
$ws.Cells.Item(18, 14).Value = 21.40931080756394

$ws.Cells.Item(19, 2).Value = 29.46527174630918
$ws.Cells.Item(19, 3).Value = 28.8037274243629
$ws.Cells.Item(19, 4).Value = 15.52054710333176
$ws.Cells.Item(19, 5).Value = 16.97617837189699
$ws.Cells.Item(19, 7).Value = 3.761269952814327
$ws.Cells.Item(19, 9).Value = 37.76802665982615
$ws.Cells.Item(19, 10).Value = 9.601461992601191
$ws.Cells.Item(19, 14).Value = 21.41180802165111

$ws.Cells.Item(20, 2).Value = 29.66543550945916
$ws.Cells.Item(20, 3).Value = 29.02070015520442
$ws.Cells.Item(20, 4).Value = 15.54434245345462
$ws.Cells.Item(20, 5).Value = 16.99986008569268
$ws.Cells.Item(20, 7).Value = 3.75901914398862
$ws.Cells.Item(20, 9).Value = 37.84765088228691
$ws.Cells.Item(20, 10).Value = 9.596820276200738
$ws.Cells.Item(20, 14).Value = 21.40070791359734

$ws.Cells.Item(21, 2).Value = 30.34032230230493
$ws.Cells.Item(21, 3).Value = 29.75022325908003
$ws.Cells.Item(21, 4).Value = 15.62900543833027
$ws.Cells.Item(21, 5).Value = 17.08464256685473
$ws.Cells.Item(21, 7).Value = 3.751668964244501
$ws.Cells.Item(21, 9).Value = 38.12545562549958
$ws.Cells.Item(21, 10).Value = 9.58332516537415
$ws.Cells.Item(21, 14).Value = 21.36583304387909

$ws.Cells.Item(22, 2).Value = 30.78240053846856
$ws.Cells.Item(22, 3).Value = 30.22657694633302
$ws.Cells.Item(22, 4).Value = 15.68782059392513
$ws.Cells.Item(22, 5).Value = 17.14391370489586
$ws.Cells.Item(22, 7).Value = 3.747022018410442
$ws.Cells.Item(22, 9).Value = 38.31457079700419
$ws.Cells.Item(22, 10).Value = 9.57607692111581
$ws.Cells.Item(22, 14).Value = 21.34485789771589

$ws.Cells.Item(23, 2).Value = 30.54641153192805
$ws.Cells.Item(23, 3).Value = 29.97242915294142
$ws.Cells.Item(23, 4).Value = 15.65611324003322
$ws.Cells.Item(23, 5).Value = 17.11192814350325
$ws.Cells.Item(23, 7).Value = 3.749487919087842
$ws.Cells.Item(23, 9).Value = 38.212952547024
$ws.Cells.Item(23, 10).Value = 9.579801256879557
$ws.Cells.Item(23, 14).Value = 21.35588551110344

$ws.Cells.Item(24, 2).Value = 29.65443260713159
$ws.Cells.Item(24, 3).Value = 29.00878093334438
$ws.Cells.Item(24, 4).Value = 15.54301792187696
$ws.Cells.Item(24, 5).Value = 16.99853991978769
$ws.Cells.Item(24, 7).Value = 3.759141945800381
$ws.Cells.Item(24, 9).Value = 37.84323937731612
$ws.Cells.Item(24, 10).Value = 9.597067291043942
$ws.Cells.Item(24, 14).Value = 21.40130840964655

$ws.Cells.Item(25, 2).Value = 28.70157155150709
$ws.Cells.Item(25, 3).Value = 27.97288974737551
$ws.Cells.Item(25, 4).Value = 15.43616368494498
$ws.Cells.Item(25, 5).Value = 16.89296384101501
$ws.Cells.Item(25, 7).Value = 3.770237465884758
$ws.Cells.Item(25, 9).Value = 37.47753344751157
$ws.Cells.Item(25, 10).Value = 9.622398831786047
$ws.Cells.Item(25, 14).Value = 21.45801831046072
